# Add three new account rows to the "Export" sheet, in the positions
# required to keep the list in its existing (descending-balance) order.
#
# Final row layout (1-based, header = row 1):
#   row 14 : 005105172 / VALDIVINO / 18000        (new - before WALTER)
#   row 15 : 004419765 / WALTER    / 17007.4       (existing, shifted down)
#   row 17 : 004911541 / TIAGO     / 11589.72      (new - before LARA)
#   row 18 : 004643737 / LARA      / 9883.81       (existing, shifted down)
#   row 22 : 005242683 / LUCAS     / 5000           (new - after LUISA)
#   row 23 : 004458604 / FABIOLA   / 4073.38        (existing, shifted down)
#
# Work from the bottom of the sheet upward so that inserting a row never
# invalidates the row numbers still to be processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) LUCAS - inserted right after LUISA (original row 19), before FABIOLA.
$ws.Rows.Item(20).Insert()
$cellA20 = $ws.Cells.Item(20, 1)
$cellA20.NumberFormat = "@"
$cellA20.Value = "005242683"
$ws.Cells.Item(20, 2).Value = "LUCAS"
$ws.Cells.Item(20, 3).Value = 5000

# 2) TIAGO - inserted right before LARA (original row 16).
$ws.Rows.Item(16).Insert()
$cellA16 = $ws.Cells.Item(16, 1)
$cellA16.NumberFormat = "@"
$cellA16.Value = "004911541"
$ws.Cells.Item(16, 2).Value = "TIAGO"
$ws.Cells.Item(16, 3).Value = 11589.72

# 3) VALDIVINO - inserted right before WALTER (original row 14).
$ws.Rows.Item(14).Insert()
$cellA14 = $ws.Cells.Item(14, 1)
$cellA14.NumberFormat = "@"
$cellA14.Value = "005105172"
$ws.Cells.Item(14, 2).Value = "VALDIVINO"
$ws.Cells.Item(14, 3).Value = 18000

Write-Output "Inserted VALDIVINO, TIAGO and LUCAS rows"
